# Loan RBI, Variable Instalments
# Insert a new column before column N ("Late") on the "Repayment Schedule" sheet
# to make room for an additional (currently blank) instalment-related column,
# shifting the existing "Late" and "Outstanding" data one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column at N; existing columns N(Late)->O, O->P, P(Outstanding)->Q
$ws.Columns("N:N").Insert()

# Match the new column's stored width to the target (10 chars, non-autofit)
$ws.Columns("N:N").ColumnWidth = 9.140625

# Restore the selected cell as recorded after the edit
$ws.Range("S8").Select() | Out-Null
